$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC (27 cell updates) -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 18519962
$ws.Range("I8").Value = 25001522
$ws.Range("K8").Value = 75004566
$ws.Range("M8").Value = -75004427
$ws.Range("H38").Value = 7832.84
$ws.Range("I38").Value = 9756.691999999999
$ws.Range("K38").Value = 29270.076
$ws.Range("M38").Value = -28898.076
$ws.Range("H110").Value = 41957.125
$ws.Range("J110").Value = 41957.125
$ws.Range("L110").Value = 41957.125
$ws.Range("N110").Value = -50137.125
$ws.Range("H132").Value = 2225676
$ws.Range("I132").Value = 1269.8462
$ws.Range("J132").Value = 16684316
$ws.Range("K132").Value = 3809.5386
$ws.Range("L132").Value = 50052948
$ws.Range("M132").Value = -1279.5386
$ws.Range("N132").Value = -50058008
$ws.Range("H136").Value = 288118.16
$ws.Range("J136").Value = 288118.16
$ws.Range("L136").Value = 288118.16
$ws.Range("N136").Value = -298318.16
$ws.Range("H137").Value = 1147867
$ws.Range("I137").Value = 3667.5
$ws.Range("K137").Value = 11002.5
$ws.Range("M137").Value = -8452.5

# ----- Sheet: ARM (47 cell updates) -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5032.6865
$ws.Range("I32").Value = 3150.4263
$ws.Range("K32").Value = 3150.4263
$ws.Range("M32").Value = -2863.4263
$ws.Range("H45").Value = 9473.6
$ws.Range("I45").Value = 12611.111
$ws.Range("K45").Value = 12611.111
$ws.Range("M45").Value = -12234.111
$ws.Range("H61").Value = 838928.75
$ws.Range("I61").Value = 5687.077
$ws.Range("K61").Value = 5687.077
$ws.Range("M61").Value = -5475.077
$ws.Range("H74").Value = 2197.926
$ws.Range("J74").Value = 3317.1428
$ws.Range("L74").Value = 3317.1428
$ws.Range("N74").Value = -5065.1428
$ws.Range("H77").Value = 2197.926
$ws.Range("J77").Value = 3317.1428
$ws.Range("L77").Value = 16585.714
$ws.Range("N77").Value = -25321.714
$ws.Range("H88").Value = 2372.6667
$ws.Range("J88").Value = 2835
$ws.Range("L88").Value = 2835
$ws.Range("N88").Value = -3647
$ws.Range("H91").Value = 2372.6667
$ws.Range("J91").Value = 2835
$ws.Range("L91").Value = 2835
$ws.Range("N91").Value = -5643
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null
$ws.Range("H122").Value = 3090858
$ws.Range("I122").Value = 4116672
$ws.Range("K122").Value = 12350016
$ws.Range("M122").Value = -12347566
$ws.Range("H132").Value = 6732.08
$ws.Range("I132").Value = 6559.1816
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 19677.5448
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -17147.5448
$ws.Range("N132").Value = -29060
$ws.Range("H136").Value = 838928.75
$ws.Range("I136").Value = 5687.077
$ws.Range("K136").Value = 17061.231
$ws.Range("M136").Value = -14511.231

# ----- Sheet: BSM (16 cell updates) -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 59631.332
$ws.Range("J2").Value = 61947
$ws.Range("L2").Value = 61947
$ws.Range("N2").Value = -62173
$ws.Range("H86").Value = 40005104
$ws.Range("J86").Value = 142858370
$ws.Range("L86").Value = 142858370
$ws.Range("N86").Value = -142860616
$ws.Range("H89").Value = 40005104
$ws.Range("J89").Value = 142858370
$ws.Range("L89").Value = 714291850
$ws.Range("N89").Value = -714303082
$ws.Range("H134").Value = 1044528.4
$ws.Range("I134").Value = 2863.5881
$ws.Range("K134").Value = 8590.764299999999
$ws.Range("M134").Value = -6055.764299999999

# ----- Sheet: CRP (15 cell updates) -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6354.963
$ws.Range("J31").Value = 7508.15
$ws.Range("L31").Value = 7508.15
$ws.Range("N31").Value = -8098.15
$ws.Range("H34").Value = 6354.963
$ws.Range("J34").Value = 7508.15
$ws.Range("L34").Value = 7508.15
$ws.Range("N34").Value = -7912.15
$ws.Range("H105").Value = 2907.8
$ws.Range("I105").Value = 2009.75
$ws.Range("J105").Value = 6500
$ws.Range("K105").Value = 2009.75
$ws.Range("L105").Value = 6500
$ws.Range("M105").Value = -262.75
$ws.Range("N105").Value = -9994

# ----- Sheet: CUL (68 cell updates) -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 634.61536
$ws.Range("I5").Value = 587.5
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 1762.5
$ws.Range("L5").Value = 3600
$ws.Range("M5").Value = -1650.5
$ws.Range("N5").Value = -3824
$ws.Range("H15").Value = 174.25
$ws.Range("J15").Value = 229.83333
$ws.Range("L15").Value = 689.49999
$ws.Range("N15").Value = -969.49999
$ws.Range("H44").Value = 394.69565
$ws.Range("J44").Value = 486.46155
$ws.Range("L44").Value = 1459.38465
$ws.Range("N44").Value = -2255.38465
$ws.Range("H45").Value = 9951
$ws.Range("J45").Value = 11883.25
$ws.Range("L45").Value = 35649.75
$ws.Range("N45").Value = -36713.75
$ws.Range("H124").Value = 2000
$ws.Range("I124").Value = 2000
$ws.Range("K124").Value = 6000
$ws.Range("M124").Value = -1090
$ws.Range("H132").Value = 2152.8
$ws.Range("I132").Value = 1283
$ws.Range("J132").Value = 2864.4546
$ws.Range("K132").Value = 11547
$ws.Range("L132").Value = 25780.0914
$ws.Range("M132").Value = -9017
$ws.Range("N132").Value = -30840.0914
$ws.Range("H133").Value = 6924.478
$ws.Range("I133").Value = 4323.4287
$ws.Range("J133").Value = 8062.4375
$ws.Range("K133").Value = 12970.2861
$ws.Range("L133").Value = 24187.3125
$ws.Range("M133").Value = -7910.286100000001
$ws.Range("N133").Value = -34307.3125
$ws.Range("H134").Value = 5620.476
$ws.Range("I134").Value = 3030
$ws.Range("K134").Value = 9090
$ws.Range("M134").Value = -4020
$ws.Range("H135").Value = 634.61536
$ws.Range("I135").Value = 587.5
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 5287.5
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -2752.5
$ws.Range("N135").Value = -15870
$ws.Range("H136").Value = 1336
$ws.Range("I136").Value = 1336
$ws.Range("K136").Value = 4008
$ws.Range("M136").Value = 1092
$ws.Range("H138").Value = 45459316
$ws.Range("I138").Value = 71432070
$ws.Range("K138").Value = 214296210
$ws.Range("M138").Value = -214291070
$ws.Range("H139").Value = 6779.1763
$ws.Range("I139").Value = 6708.6665
$ws.Range("K139").Value = 20125.9995
$ws.Range("M139").Value = -14985.9995
$ws.Range("H140").Value = 65273.25
$ws.Range("I140").Value = 92215.63
$ws.Range("K140").Value = 276646.89
$ws.Range("M140").Value = -271466.89
$ws.Range("H141").Value = 333337060
$ws.Range("J141").Value = 6999
$ws.Range("L141").Value = 20997
$ws.Range("N141").Value = -31357

# ----- Sheet: GSM (20 cell updates) -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 23333
$ws.Range("H55").Value = 13000
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 14500
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 14500
$ws.Range("M55").Value = -9673
$ws.Range("N55").Value = -15154
$ws.Range("H126").Value = 3284.5625
$ws.Range("I126").Value = 2428.7778
$ws.Range("K126").Value = 7286.3334
$ws.Range("M126").Value = -4816.3334
$ws.Range("H128").Value = 56695
$ws.Range("J128").Value = 56695
$ws.Range("L128").Value = 56695
$ws.Range("N128").Value = -66655
$ws.Range("H132").Value = 43480984
$ws.Range("J132").Value = 8992.6
$ws.Range("L132").Value = 26977.8
$ws.Range("N132").Value = -32037.8

# ----- Sheet: LTW (37 cell updates) -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3597.3684
$ws.Range("I7").Value = 3436.6667
$ws.Range("J7").Value = 4200
$ws.Range("K7").Value = 3436.6667
$ws.Range("L7").Value = 4200
$ws.Range("M7").Value = -3324.6667
$ws.Range("N7").Value = -4424
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = $null
$ws.Range("H61").Value = 3847.7917
$ws.Range("I61").Value = 2138.0588
$ws.Range("K61").Value = 2138.0588
$ws.Range("M61").Value = -1936.0588
$ws.Range("H93").Value = 4424.294
$ws.Range("I93").Value = 1776.5
$ws.Range("K93").Value = 1776.5
$ws.Range("M93").Value = -528.5
$ws.Range("H113").Value = 3847.7917
$ws.Range("I113").Value = 2138.0588
$ws.Range("K113").Value = 2138.0588
$ws.Range("M113").Value = 31.94120000000021
$ws.Range("H126").Value = 3597.3684
$ws.Range("I126").Value = 3436.6667
$ws.Range("J126").Value = 4200
$ws.Range("K126").Value = 10310.0001
$ws.Range("L126").Value = 12600
$ws.Range("M126").Value = -7840.000100000001
$ws.Range("N126").Value = -17540
$ws.Range("H132").Value = 5500
$ws.Range("I132").Value = 4500
$ws.Range("J132").Value = 5750
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 17250
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -22310

# ----- Sheet: WVR (8 cell updates) -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 650
$ws.Range("J13").Value = 475
$ws.Range("L13").Value = 475
$ws.Range("N13").Value = -755
$ws.Range("H126").Value = 9806553
$ws.Range("J126").Value = 18521172
$ws.Range("L126").Value = 55563516
$ws.Range("N126").Value = -55568456
